$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and shift rankings to insert OKB at row 8).
# Only cells whose content actually changes are touched; row A (rank index) and
# header row 1 are untouched, matching the source diff.

$ws.Range("D2").Value = "25.770.09"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.751.12"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'236.91"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5062"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'40.39"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2605"
$ws.Range("E9").Value = "  +7.80%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06197"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.751.43"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.06938"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'15.46"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6054"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'78.51"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.459"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "25.789.24"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'11.65"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000006733"
$ws.Range("E21").Value = "  +6.24%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "1.974.89"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.051"
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'8.193"
$ws.Range("E24").Value = "  +4.70%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'5.153"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'137.99"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.460"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.10"
$ws.Range("E28").Value = "  +4.88%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.804"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'102.02"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08256"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'3.706"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.394"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04352"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "'0.9999"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.650"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.6004"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.701"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.955"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01547"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'103.09"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.7536"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3803"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'4.858"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05485"
$ws.Range("E47").Value = "  +7.74%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1076"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'5.933"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'30.12"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.52%  "
